# Insert a new weekly data row above row 121 (Apio / Feria Lagunitas de Puerto
# Montt), which pushes every existing row from 121..235 down by one (new
# 122..236) and grows the used range from A1:R235 to A1:R236.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(121).Insert()

# Populate the newly inserted row with the new weekly reading. Every column
# except Fecha (D) and Volumen (J) mirrors the entry that used to sit in row
# 121 (same market/region/product/quality/prices/unit/origin combination).
$ws.Cells(121, 1).Value = 4
$ws.Cells(121, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells(121, 3).Value = "Los Lagos"
$ws.Cells(121, 4).Value = 44669
$ws.Cells(121, 5).Value = 10
$ws.Cells(121, 6).Value = 100112017
$ws.Cells(121, 7).Value = "Apio"
$ws.Cells(121, 8).Value = "Americana (o)"
$ws.Cells(121, 9).Value = "Primera"
$ws.Cells(121, 10).Value = 10
$ws.Cells(121, 11).Value = 12000
$ws.Cells(121, 12).Value = 12000
$ws.Cells(121, 13).Value = 12000
$ws.Cells(121, 14).Value = "`$/docena de matas"
$ws.Cells(121, 15).Value = "Región de Coquimbo"
$ws.Cells(121, 16).Value = 2000
$ws.Cells(121, 17).Value = 6
$ws.Cells(121, 18).Value = "Hortaliza"
